# Apply updated cryptocurrency price/volume data to sheet1 (cryptos.xlsx)
# All target cells are plain text cells (t="inlineStr" in the original file),
# so we force text storage via NumberFormat "@" (to stop Excel auto-parsing
# numeric-looking strings such as "1.002" or "30.691.65" into numbers), then
# reset the style back to "Normal" so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "30.691.65"
Set-TextValue 2 5 "  +1.93%  "

# Row 3
Set-TextValue 3 4 "1.894.88"
Set-TextValue 3 5 "  +1.09%  "

# Row 4
Set-TextValue 4 4 "1.002"
Set-TextValue 4 5 "  +0.19%  "

# Row 5
Set-TextValue 5 4 "242.15"
Set-TextValue 5 5 "  -0.02%  "

# Row 6
Set-TextValue 6 4 "1.001"
Set-TextValue 6 5 "  +0.06%  "

# Row 7
Set-TextValue 7 4 "0.4921"
Set-TextValue 7 5 "  +0.70%  "

# Row 8
Set-TextValue 8 5 "  +1.33%  "

# Row 9
Set-TextValue 9 4 "0.06741"
Set-TextValue 9 5 "  +2.47%  "

# Row 10
Set-TextValue 10 4 "1.897.47"
Set-TextValue 10 5 "  +1.19%  "

# Row 11
Set-TextValue 11 4 "17.20"
Set-TextValue 11 5 "  +5.07%  "

# Row 12
Set-TextValue 12 4 "0.07248"
Set-TextValue 12 5 "  +0.72%  "

# Row 13
Set-TextValue 13 4 "90.65"
Set-TextValue 13 5 "  +5.71%  "

# Row 14
Set-TextValue 14 4 "0.6756"
Set-TextValue 14 5 "  +1.74%  "

# Row 15
Set-TextValue 15 4 "5.021"
Set-TextValue 15 5 "  +2.40%  "

# Row 16
Set-TextValue 16 4 "30.684.89"
Set-TextValue 16 5 "  +2.03%  "

# Row 17
Set-TextValue 17 4 "0.000007970"
Set-TextValue 17 5 "  +2.42%  "

# Row 18
Set-TextValue 18 5 "  +0.08%  "

# Row 19
Set-TextValue 19 4 "13.07"
Set-TextValue 19 5 "  +2.77%  "

# Row 20
Set-TextValue 20 4 "2.142.48"
Set-TextValue 20 5 "  +1.12%  "

# Row 21
Set-TextValue 21 4 "1.002"
Set-TextValue 21 5 "  +0.16%  "

# Row 22
Set-TextValue 22 4 "4.803"
Set-TextValue 22 5 "  +1.24%  "

# Row 23
Set-TextValue 23 4 "188.86"
Set-TextValue 23 5 "  +31.99%  "

# Row 24
Set-TextValue 24 4 "6.079"
Set-TextValue 24 5 "  +4.38%  "

# Row 25
Set-TextValue 25 4 "9.360"
Set-TextValue 25 5 "  +2.05%  "

# Row 26
Set-TextValue 26 4 "156.53"
Set-TextValue 26 5 "  +2.43%  "

# Row 27
Set-TextValue 27 4 "18.91"
Set-TextValue 27 5 "  +11.74%  "

# Row 28
Set-TextValue 28 4 "1.892"
Set-TextValue 28 5 "  +1.00%  "

# Row 29
Set-TextValue 29 5 "  +0.39%  "

# Row 30
Set-TextValue 30 4 "4.294"
Set-TextValue 30 5 "  +2.40%  "

# Row 31
Set-TextValue 31 4 "0.09067"
Set-TextValue 31 5 "  +3.48%  "

# Row 32
Set-TextValue 32 4 "3.993"
Set-TextValue 32 5 "  +0.24%  "

# Row 33
Set-TextValue 33 4 "0.05220"
Set-TextValue 33 5 "  +1.91%  "

# Row 34
Set-TextValue 34 4 "0.7357"
Set-TextValue 34 5 "  +3.54%  "

# Row 35
Set-TextValue 35 4 "1.106"
Set-TextValue 35 5 "  +0.17%  "

# Row 36
Set-TextValue 36 4 "2.761"
Set-TextValue 36 5 "  +3.45%  "

# Row 37
Set-TextValue 37 4 "0.01829"
Set-TextValue 37 5 "  -0.46%  "

# Row 38
Set-TextValue 38 4 "2.680"
Set-TextValue 38 5 "  +0.08%  "

# Row 39
Set-TextValue 39 4 "2.122"
Set-TextValue 39 5 "  -0.16%  "

# Row 40
Set-TextValue 40 4 "0.9303"
Set-TextValue 40 5 "  +0.83%  "

# Row 41
Set-TextValue 41 4 "0.4390"
Set-TextValue 41 5 "  +4.26%  "

# Row 42
Set-TextValue 42 4 "105.02"
Set-TextValue 42 5 "  +1.11%  "

# Row 44
Set-TextValue 44 4 "5.725"
Set-TextValue 44 5 "  -0.54%  "

# Row 45
Set-TextValue 45 2 "Algorand"
Set-TextValue 45 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 45 4 "0.1350"
Set-TextValue 45 5 "  +5.62%  "

# Row 46
Set-TextValue 46 2 "Aptos"
Set-TextValue 46 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 46 4 "7.519"
Set-TextValue 46 5 "  +1.67%  "

# Row 47
Set-TextValue 47 4 "0.05866"
Set-TextValue 47 5 "  +2.68%  "

# Row 48
Set-TextValue 48 4 "8.751"
Set-TextValue 48 5 "  +6.75%  "

# Row 49
Set-TextValue 49 4 "0.3930"
Set-TextValue 49 5 "  +4.91%  "

# Row 50
Set-TextValue 50 4 "33.61"
Set-TextValue 50 5 "  +2.58%  "

# Row 51
Set-TextValue 51 4 "1.415"
Set-TextValue 51 5 "  +6.21%  "
